$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the generated Appointment ID for the row that previously held the
# placeholder "test" value so the appointment list has a unique, real ID
# (able to generate unique appointment ID).
$ws.Range("A5").Value = "A004"

# Move the active selection to match the cursor position left by the edit.
$ws.Range("B9").Select()
